$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$question = "I have 20000 modifiers added ty log, why I can't I add anymore?"
$model = "llama3.2:latest"

$responseA = "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file. The maximum number of layouts allowed is 100 tables in a log."
$responseB = "You cannot add more than 20000 modifiers per plot because of the limit on the number of layouts per ODF file."

$rows = @(
    @{ Row = 274; Response = $responseA },
    @{ Row = 275; Response = $responseA },
    @{ Row = 276; Response = $responseA },
    @{ Row = 277; Response = $responseB },
    @{ Row = 278; Response = $responseB }
)

foreach ($entry in $rows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $question
    $ws.Cells.Item($r, 2).Value = $model
    $ws.Cells.Item($r, 3).Value = $entry.Response
}
